$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 updates
$ws.Range("G8").Value = 1.57
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 2.1
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 12
$ws.Range("W8").Value = 7
$ws.Range("X8").Value = 7.5
$ws.Range("AC8").Value = 12
$ws.Range("AD8").Value = 8
$ws.Range("AG8").Value = 700
$ws.Range("AH8").Value = 15
$ws.Range("AN8").Value = 3.6
$ws.Range("AO8").Value = 8
$ws.Range("AQ8").Value = 23
$ws.Range("AR8").Value = 41
$ws.Range("AX8").Value = 7
$ws.Range("BC8").Value = 450

# Row 9 updates
$ws.Range("G9").Value = 8
$ws.Range("H9").Value = 4.75
$ws.Range("I9").Value = 1.36
$ws.Range("J9").Value = 6.5
$ws.Range("K9").Value = 2.6
$ws.Range("L9").Value = 1.8
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 10.5
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.38
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("W9").Value = 23
$ws.Range("Y9").Value = 23
$ws.Range("Z9").Value = 81
$ws.Range("AD9").Value = 9.5
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 9
$ws.Range("AJ9").Value = 9
$ws.Range("AK9").Value = 9.5
$ws.Range("AN9").Value = 8.5
$ws.Range("AO9").Value = 34
$ws.Range("AP9").Value = 34
$ws.Range("AQ9").Value = 126
$ws.Range("AR9").Value = 126
$ws.Range("AS9").Value = 400
$ws.Range("AU9").Value = 8.5
$ws.Range("AX9").Value = 3.6
$ws.Range("AY9").Value = 6.5
$ws.Range("BA9").Value = 17
